$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.648.57"
$ws.Range("E2").Value = "  -2.21%  "
$ws.Range("D3").Value = "1.799.48"
$ws.Range("E3").Value = "  -1.79%  "
$ws.Range("E4").Value = "  -0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "231.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5910"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.99%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2770"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06808"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.64%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07517"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.97%  "
$ws.Range("D12").Value = "1.788.08"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  -0.92%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6212"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.25%  "
$ws.Range("D15").Value = "2.044.57"
$ws.Range("E15").Value = "  -1.74%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009130"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.19%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "75.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.83%  "
$ws.Range("D18").Value = "28.623.80"
$ws.Range("E18").Value = "  -2.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.461"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.76%  "
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "210.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -6.90%  "
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.824"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("E24").Value = "  -0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.851"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.24%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1269"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.43"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.418"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06192"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.31%  "
$ws.Range("E31").Value = "  -1.67%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.802"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.776"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.69%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.730"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.93%  "
$ws.Range("E35").Value = "  -6.08%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6419"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.38%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.498"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.710"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.03%  "
$ws.Range("E39").Value = "  -0.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01698"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.74%  "
$ws.Range("D41").Value = "1.146.10"
$ws.Range("E41").Value = "  -6.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8819"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.008"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.15"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.73%  "
$ws.Range("D45").Value = "1.947.95"
$ws.Range("E45").Value = "  -2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.30"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.81%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000112"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05467"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.90%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.320"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.90%  "
$ws.Range("E51").Value = "  -1.99%  "
